$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "What is the required admin features?" (row 7) proposal text is updated
# to mention the new "search" capability alongside add/delete.
$ws.Range("C7").Value = "Admin can add ,search and delete users."

# Leave the cursor where the edit was made, matching the author's final
# on-screen selection.
$ws.Range("C7").Select()
